$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data values
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.459.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4588"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3829"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07924"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9669"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.829.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.879"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.052"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06657"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001031"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.465.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.338"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.307"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.082.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.062"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.245"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9496"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09286"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.566"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.234"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05931"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02196"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.155"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.007"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5789"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1840"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.276"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5487"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.867"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06641"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "109.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.040"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
